## Applies the commit:
##  1) Re-points the three data tables (slides 14, 15, 16) from the
##     bespoke "Table_0" style onto the built-in PowerPoint table style
##     {D592A325-DB20-4E6B-89C2-183F40E030B8}.
##  2) Swaps the colour scheme carried by the deck's theme (theme1.xml,
##     the one bound to the slide master) so it now holds the stock
##     "Office" palette instead of the "Integral"/"Red Violet" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style re-application -----------------------------------
$newStyleId = "{D592A325-DB20-4E6B-89C2-183F40E030B8}"
$tableSlides = 14, 15, 16
foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour-scheme swap --------------------------------------
# theme1.xml / theme2.xml already share an identical font scheme and
# format scheme; only the 12 colour-scheme slots (and the cosmetic
# "name" attributes, which aren't exposed through the object model)
# differ between the "Integral" theme and the "Office Theme".
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
